$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.305.44'
$ws.Range('E2').Value = '  +1.49%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.808.34'
$ws.Range('E3').Value = '  +3.58%  '

$ws.Range('E4').Value = '  -0.33%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '338.25'
$ws.Range('E5').Value = '  +1.20%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9997'
$ws.Range('E6').Value = '  -0.06%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4666'
$ws.Range('E7').Value = '  +21.60%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3832'
$ws.Range('E8').Value = '  +13.49%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.54'
$ws.Range('E9').Value = '  -0.79%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.160'
$ws.Range('E10').Value = '  +4.71%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07627'
$ws.Range('E11').Value = '  +6.45%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.53'
$ws.Range('E12').Value = '  +1.08%  '

$ws.Range('E13').Value = '  -0.32%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.360'
$ws.Range('E14').Value = '  +3.84%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.457'
$ws.Range('E15').Value = '  +5.49%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.808.22'
$ws.Range('E16').Value = '  +3.25%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001097'
$ws.Range('E17').Value = '  +4.19%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06715'
$ws.Range('E18').Value = '  +1.76%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '81.80'
$ws.Range('E19').Value = '  +3.80%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9994'
$ws.Range('E20').Value = '  -0.21%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.59'
$ws.Range('E21').Value = '  +5.59%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.428'
$ws.Range('E22').Value = '  +4.57%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.299.77'
$ws.Range('E23').Value = '  +1.34%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.91'
$ws.Range('E24').Value = '  +2.93%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.423'
$ws.Range('E25').Value = '  +1.00%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '20.75'
$ws.Range('E26').Value = '  +5.12%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '153.30'
$ws.Range('E27').Value = '  +0.16%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.382'
$ws.Range('E28').Value = '  +4.37%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.013.87'
$ws.Range('E29').Value = '  +3.25%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.44'
$ws.Range('E30').Value = '  +2.36%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.260'
$ws.Range('E31').Value = '  -1.38%  '

$ws.Range('E32').Value = '  +0.27%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09561'
$ws.Range('E33').Value = '  +8.99%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.869'
$ws.Range('E34').Value = '  +1.59%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.2285'
$ws.Range('E35').Value = '  +9.91%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06392'
$ws.Range('E36').Value = '  +5.12%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '12.12'
$ws.Range('E37').Value = '  +0.31%  '

$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02359'
$ws.Range('E38').Value = '  +4.02%  '

$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.282'
$ws.Range('E39').Value = '  +3.56%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6648'
$ws.Range('E40').Value = '  +2.42%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.241'
$ws.Range('E41').Value = '  +3.50%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.494'
$ws.Range('E42').Value = '  -3.00%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.342'
$ws.Range('E43').Value = '  +5.50%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.25'
$ws.Range('E44').Value = '  +4.68%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9994'
$ws.Range('E45').Value = '  -0.10%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6152'
$ws.Range('E46').Value = '  +2.82%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.857'
$ws.Range('E47').Value = '  +0.72%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '131.02'
$ws.Range('E48').Value = '  +3.63%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.043'
$ws.Range('E49').Value = '  +3.02%  '

$ws.Range('E50').Value = '  +1.68%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07150'
$ws.Range('E51').Value = '  +2.53%  '
